# Realestate Update resale numbers 2023-06-22 18:43
# Append a new data row (row 68) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D are text-like (date/time/weekday/week-number) in this sheet,
# so prefix numeric-looking text with a leading apostrophe to keep them
# stored as text instead of being auto-converted to a date/number by Excel.
$ws.Cells.Item(68, 1).Value = "'2023-06-22"
$ws.Cells.Item(68, 2).Value = "18:43:49"
$ws.Cells.Item(68, 3).Value = "Thursday"
$ws.Cells.Item(68, 4).Value = "'25"
$ws.Cells.Item(68, 5).Value = 122460
$ws.Cells.Item(68, 6).Value = 133784
$ws.Cells.Item(68, 7).Value = 162573
$ws.Cells.Item(68, 8).Value = 133608
$ws.Cells.Item(68, 9).Value = 177361
$ws.Cells.Item(68, 10).Value = 115102
$ws.Cells.Item(68, 11).Value = 202417
$ws.Cells.Item(68, 12).Value = 225843
$ws.Cells.Item(68, 13).Value = 175498
$ws.Cells.Item(68, 14).Value = 104026
$ws.Cells.Item(68, 15).Value = 39396
$ws.Cells.Item(68, 16).Value = 33876
$ws.Cells.Item(68, 17).Value = 51934
$ws.Cells.Item(68, 18).Value = -1
$ws.Cells.Item(68, 19).Value = 36156
$ws.Cells.Item(68, 20).Value = -1
